$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.ClearFormats()
}

Set-TextCell 'D2' '43.052.69'
Set-TextCell 'E2' '  -6.92%  '

Set-TextCell 'D3' '2.549.48'
Set-TextCell 'E3' '  -2.58%  '

Set-TextCell 'E4' '  +0.02%  '

Set-TextCell 'D5' '298.50'
Set-TextCell 'E5' '  -4.11%  '

Set-TextCell 'D6' '92.66'
Set-TextCell 'E6' '  -7.03%  '

Set-TextCell 'D7' '0.575'
Set-TextCell 'E7' '  -3.80%  '

Set-TextCell 'E8' '  -0.01%  '

Set-TextCell 'D9' '0.552'
Set-TextCell 'E9' '  -5.47%  '

Set-TextCell 'D10' '36.01'
Set-TextCell 'E10' '  -7.86%  '

Set-TextCell 'E11' '  -4.39%  '

Set-TextCell 'D12' '7.74'
Set-TextCell 'E12' '  -4.86%  '

Set-TextCell 'E13' '  +1.29%  '

Set-TextCell 'D14' '2.939.23'
Set-TextCell 'E14' '  -2.54%  '

Set-TextCell 'D15' '2.549.61'
Set-TextCell 'E15' '  -2.63%  '

Set-TextCell 'D16' '0.871'
Set-TextCell 'E16' '  -5.24%  '

Set-TextCell 'D17' '14.21'
Set-TextCell 'E17' '  -4.70%  '

Set-TextCell 'D18' '43.089.51'
Set-TextCell 'E18' '  -7.36%  '

Set-TextCell 'D19' '12.71'
Set-TextCell 'E19' '  -1.16%  '

Set-TextCell 'B20' 'ShibaInu'
Set-TextCell 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D20' '0.0₃0980'
Set-TextCell 'E20' '  -4.04%  '

Set-TextCell 'B21' 'Uniswap'
Set-TextCell 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D21' '6.66'
Set-TextCell 'E21' '  -2.12%  '

Set-TextCell 'D22' '71.78'
Set-TextCell 'E22' '  -2.20%  '

Set-TextCell 'D23' '260.88'
Set-TextCell 'E23' '  -10.93%  '

Set-TextCell 'E24' '  -4.81%  '

Set-TextCell 'D25' '2.15'
Set-TextCell 'E25' '  -4.98%  '

Set-TextCell 'D26' '29.58'
Set-TextCell 'E26' '  -0.51%  '

Set-TextCell 'E27' '  +0.03%  '

Set-TextCell 'D28' '10.08'
Set-TextCell 'E28' '  -7.41%  '

Set-TextCell 'D29' '36.84'
Set-TextCell 'E29' '  -6.14%  '

Set-TextCell 'E30' '  -4.04%  '

Set-TextCell 'D31' '5.97'
Set-TextCell 'E31' '  -4.89%  '

Set-TextCell 'D32' '155.64'
Set-TextCell 'E32' '  -1.81%  '

Set-TextCell 'D33' '2.17'

Set-TextCell 'D34' '3.40'
Set-TextCell 'E34' '  -5.82%  '

Set-TextCell 'D35' '2.73'

Set-TextCell 'D36' '0.0799'
Set-TextCell 'E36' '  -5.48%  '

Set-TextCell 'E37' '  -5.79%  '

Set-TextCell 'E38' '  -3.28%  '

Set-TextCell 'D39' '23.57'
Set-TextCell 'E39' '  +9.27%  '

Set-TextCell 'D40' '16.48'
Set-TextCell 'E40' '  +3.90%  '

Set-TextCell 'D41' '3.47'
Set-TextCell 'E41' '  -3.56%  '

Set-TextCell 'E42' '  -5.57%  '

Set-TextCell 'D43' '3.87'
Set-TextCell 'E43' '  -4.11%  '

Set-TextCell 'D44' '2.075.31'
Set-TextCell 'E44' '  -2.22%  '

Set-TextCell 'D45' '0.998'
Set-TextCell 'E45' '  -0.04%  '

Set-TextCell 'D46' '85.84'
Set-TextCell 'E46' '  -11.65%  '

Set-TextCell 'E47' '  +3.06%  '

Set-TextCell 'B48' 'RocketPoolETH'
Set-TextCell 'C48' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 'D48' '2.796.48'
Set-TextCell 'E48' '  -2.69%  '

Set-TextCell 'B49' 'Stacks'
Set-TextCell 'C49' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 'D49' '1.71'
Set-TextCell 'E49' '  -2.26%  '

Set-TextCell 'D50' '8.77'
Set-TextCell 'E50' '  -8.39%  '

Set-TextCell 'D51' '104.35'
Set-TextCell 'E51' '  -4.99%  '
